$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "White - Asian"
$ws.Cells.Item(2, 2).Value = 2.53
$ws.Cells.Item(2, 3).Value = 0.604
$ws.Cells.Item(2, 4).Value = 1669
$ws.Cells.Item(2, 5).Value = 4.187
$ws.Cells.Item(2, 6).Value = "'" + "2.8701e-04"
$ws.Cells.Item(2, 7).Value = "'" + "TRUE"
$ws.Cells.Item(2, 8).Value = "AA_C"

$ws.Cells.Item(3, 1).Value = "White - Black"
$ws.Cells.Item(3, 2).Value = -0.62
$ws.Cells.Item(3, 3).Value = 0.296
$ws.Cells.Item(3, 4).Value = 1669
$ws.Cells.Item(3, 5).Value = -2.092
$ws.Cells.Item(3, 6).Value = "'" + "2.2380e-01"
$ws.Cells.Item(3, 7).Value = "'" + "FALSE"
$ws.Cells.Item(3, 8).Value = "AA_C"

$ws.Cells.Item(4, 1).Value = "White - LatinX"
$ws.Cells.Item(4, 2).Value = 0.096
$ws.Cells.Item(4, 3).Value = 0.593
$ws.Cells.Item(4, 4).Value = 1669
$ws.Cells.Item(4, 5).Value = 0.162
$ws.Cells.Item(4, 6).Value = "'" + "9.9985e-01"
$ws.Cells.Item(4, 7).Value = "'" + "FALSE"
$ws.Cells.Item(4, 8).Value = "AA_C"

$ws.Cells.Item(5, 1).Value = "White - Other"
$ws.Cells.Item(5, 2).Value = 0.135
$ws.Cells.Item(5, 3).Value = 1.039
$ws.Cells.Item(5, 4).Value = 1669
$ws.Cells.Item(5, 5).Value = 0.13
$ws.Cells.Item(5, 6).Value = "'" + "9.9994e-01"
$ws.Cells.Item(5, 7).Value = "'" + "FALSE"
$ws.Cells.Item(5, 8).Value = "AA_C"

$ws.Cells.Item(6, 1).Value = "Asian - Black"
$ws.Cells.Item(6, 2).Value = -3.15
$ws.Cells.Item(6, 3).Value = 0.633
$ws.Cells.Item(6, 4).Value = 1669
$ws.Cells.Item(6, 5).Value = -4.974
$ws.Cells.Item(6, 6).Value = "'" + "7.1779e-06"
$ws.Cells.Item(6, 7).Value = "'" + "TRUE"
$ws.Cells.Item(6, 8).Value = "AA_C"

$ws.Cells.Item(7, 1).Value = "Asian - LatinX"
$ws.Cells.Item(7, 2).Value = -2.434
$ws.Cells.Item(7, 3).Value = 0.813
$ws.Cells.Item(7, 4).Value = 1669
$ws.Cells.Item(7, 5).Value = -2.993
$ws.Cells.Item(7, 6).Value = "'" + "2.3429e-02"
$ws.Cells.Item(7, 7).Value = "'" + "TRUE"
$ws.Cells.Item(7, 8).Value = "AA_C"

$ws.Cells.Item(8, 1).Value = "Asian - Other"
$ws.Cells.Item(8, 2).Value = -2.396
$ws.Cells.Item(8, 3).Value = 1.18
$ws.Cells.Item(8, 4).Value = 1669
$ws.Cells.Item(8, 5).Value = -2.03
$ws.Cells.Item(8, 6).Value = "'" + "2.5198e-01"
$ws.Cells.Item(8, 7).Value = "'" + "FALSE"
$ws.Cells.Item(8, 8).Value = "AA_C"

$ws.Cells.Item(9, 1).Value = "Black - LatinX"
$ws.Cells.Item(9, 2).Value = 0.716
$ws.Cells.Item(9, 3).Value = 0.622
$ws.Cells.Item(9, 4).Value = 1669
$ws.Cells.Item(9, 5).Value = 1.15
$ws.Cells.Item(9, 6).Value = "'" + "7.7952e-01"
$ws.Cells.Item(9, 7).Value = "'" + "FALSE"
$ws.Cells.Item(9, 8).Value = "AA_C"

$ws.Cells.Item(10, 1).Value = "Black - Other"
$ws.Cells.Item(10, 2).Value = 0.754
$ws.Cells.Item(10, 3).Value = 1.055
$ws.Cells.Item(10, 4).Value = 1669
$ws.Cells.Item(10, 5).Value = 0.715
$ws.Cells.Item(10, 6).Value = "'" + "9.5312e-01"
$ws.Cells.Item(10, 7).Value = "'" + "FALSE"
$ws.Cells.Item(10, 8).Value = "AA_C"

$ws.Cells.Item(11, 1).Value = "LatinX - Other"
$ws.Cells.Item(11, 2).Value = 0.039
$ws.Cells.Item(11, 3).Value = 1.174
$ws.Cells.Item(11, 4).Value = 1669
$ws.Cells.Item(11, 5).Value = 0.033
$ws.Cells.Item(11, 6).Value = "'" + "1.0000e+00"
$ws.Cells.Item(11, 7).Value = "'" + "FALSE"
$ws.Cells.Item(11, 8).Value = "AA_C"

$ws.Cells.Item(12, 1).Value = "White - Asian"
$ws.Cells.Item(12, 2).Value = -0.723
$ws.Cells.Item(12, 3).Value = 1.065
$ws.Cells.Item(12, 4).Value = 1669
$ws.Cells.Item(12, 5).Value = -0.679
$ws.Cells.Item(12, 6).Value = "'" + "9.6108e-01"
$ws.Cells.Item(12, 7).Value = "'" + "FALSE"
$ws.Cells.Item(12, 8).Value = "BiW_L"

$ws.Cells.Item(13, 1).Value = "White - Black"
$ws.Cells.Item(13, 2).Value = -4.782
$ws.Cells.Item(13, 3).Value = 0.522
$ws.Cells.Item(13, 4).Value = 1669
$ws.Cells.Item(13, 5).Value = -9.158
$ws.Cells.Item(13, 6).Value = "'" + "1.6549e-12"
$ws.Cells.Item(13, 7).Value = "'" + "TRUE"
$ws.Cells.Item(13, 8).Value = "BiW_L"

$ws.Cells.Item(14, 1).Value = "White - LatinX"
$ws.Cells.Item(14, 2).Value = -1.57
$ws.Cells.Item(14, 3).Value = 1.045
$ws.Cells.Item(14, 4).Value = 1669
$ws.Cells.Item(14, 5).Value = -1.503
$ws.Cells.Item(14, 6).Value = "'" + "5.6073e-01"
$ws.Cells.Item(14, 7).Value = "'" + "FALSE"
$ws.Cells.Item(14, 8).Value = "BiW_L"

$ws.Cells.Item(15, 1).Value = "White - Other"
$ws.Cells.Item(15, 2).Value = -0.878
$ws.Cells.Item(15, 3).Value = 1.831
$ws.Cells.Item(15, 4).Value = 1669
$ws.Cells.Item(15, 5).Value = -0.48
$ws.Cells.Item(15, 6).Value = "'" + "9.8921e-01"
$ws.Cells.Item(15, 7).Value = "'" + "FALSE"
$ws.Cells.Item(15, 8).Value = "BiW_L"

$ws.Cells.Item(16, 1).Value = "Asian - Black"
$ws.Cells.Item(16, 2).Value = -4.059
$ws.Cells.Item(16, 3).Value = 1.116
$ws.Cells.Item(16, 4).Value = 1669
$ws.Cells.Item(16, 5).Value = -3.637
$ws.Cells.Item(16, 6).Value = "'" + "2.6345e-03"
$ws.Cells.Item(16, 7).Value = "'" + "TRUE"
$ws.Cells.Item(16, 8).Value = "BiW_L"

$ws.Cells.Item(17, 1).Value = "Asian - LatinX"
$ws.Cells.Item(17, 2).Value = -0.847
$ws.Cells.Item(17, 3).Value = 1.433
$ws.Cells.Item(17, 4).Value = 1669
$ws.Cells.Item(17, 5).Value = -0.591
$ws.Cells.Item(17, 6).Value = "'" + "9.7642e-01"
$ws.Cells.Item(17, 7).Value = "'" + "FALSE"
$ws.Cells.Item(17, 8).Value = "BiW_L"

$ws.Cells.Item(18, 1).Value = "Asian - Other"
$ws.Cells.Item(18, 2).Value = -0.155
$ws.Cells.Item(18, 3).Value = 2.08
$ws.Cells.Item(18, 4).Value = 1669
$ws.Cells.Item(18, 5).Value = -0.075
$ws.Cells.Item(18, 6).Value = "'" + "9.9999e-01"
$ws.Cells.Item(18, 7).Value = "'" + "FALSE"
$ws.Cells.Item(18, 8).Value = "BiW_L"

$ws.Cells.Item(19, 1).Value = "Black - LatinX"
$ws.Cells.Item(19, 2).Value = 3.212
$ws.Cells.Item(19, 3).Value = 1.097
$ws.Cells.Item(19, 4).Value = 1669
$ws.Cells.Item(19, 5).Value = 2.929
$ws.Cells.Item(19, 6).Value = "'" + "2.8430e-02"
$ws.Cells.Item(19, 7).Value = "'" + "TRUE"
$ws.Cells.Item(19, 8).Value = "BiW_L"

$ws.Cells.Item(20, 1).Value = "Black - Other"
$ws.Cells.Item(20, 2).Value = 3.904
$ws.Cells.Item(20, 3).Value = 1.86
$ws.Cells.Item(20, 4).Value = 1669
$ws.Cells.Item(20, 5).Value = 2.099
$ws.Cells.Item(20, 6).Value = "'" + "2.2081e-01"
$ws.Cells.Item(20, 7).Value = "'" + "FALSE"
$ws.Cells.Item(20, 8).Value = "BiW_L"

$ws.Cells.Item(21, 1).Value = "LatinX - Other"
$ws.Cells.Item(21, 2).Value = 0.692
$ws.Cells.Item(21, 3).Value = 2.068
$ws.Cells.Item(21, 4).Value = 1669
$ws.Cells.Item(21, 5).Value = 0.335
$ws.Cells.Item(21, 6).Value = "'" + "9.9730e-01"
$ws.Cells.Item(21, 7).Value = "'" + "FALSE"
$ws.Cells.Item(21, 8).Value = "BiW_L"

$ws.Cells.Item(22, 1).Value = "White - Asian"
$ws.Cells.Item(22, 2).Value = 5.34
$ws.Cells.Item(22, 3).Value = 1.377
$ws.Cells.Item(22, 4).Value = 1669
$ws.Cells.Item(22, 5).Value = 3.878
$ws.Cells.Item(22, 6).Value = "'" + "1.0340e-03"
$ws.Cells.Item(22, 7).Value = "'" + "TRUE"
$ws.Cells.Item(22, 8).Value = "BiW_C"

$ws.Cells.Item(23, 1).Value = "White - Black"
$ws.Cells.Item(23, 2).Value = -1.35
$ws.Cells.Item(23, 3).Value = 0.675
$ws.Cells.Item(23, 4).Value = 1669
$ws.Cells.Item(23, 5).Value = -2.0
$ws.Cells.Item(23, 6).Value = "'" + "2.6613e-01"
$ws.Cells.Item(23, 7).Value = "'" + "FALSE"
$ws.Cells.Item(23, 8).Value = "BiW_C"

$ws.Cells.Item(24, 1).Value = "White - LatinX"
$ws.Cells.Item(24, 2).Value = -0.295
$ws.Cells.Item(24, 3).Value = 1.351
$ws.Cells.Item(24, 4).Value = 1669
$ws.Cells.Item(24, 5).Value = -0.219
$ws.Cells.Item(24, 6).Value = "'" + "9.9949e-01"
$ws.Cells.Item(24, 7).Value = "'" + "FALSE"
$ws.Cells.Item(24, 8).Value = "BiW_C"

$ws.Cells.Item(25, 1).Value = "White - Other"
$ws.Cells.Item(25, 2).Value = 0.662
$ws.Cells.Item(25, 3).Value = 2.367
$ws.Cells.Item(25, 4).Value = 1669
$ws.Cells.Item(25, 5).Value = 0.28
$ws.Cells.Item(25, 6).Value = "'" + "9.9866e-01"
$ws.Cells.Item(25, 7).Value = "'" + "FALSE"
$ws.Cells.Item(25, 8).Value = "BiW_C"

$ws.Cells.Item(26, 1).Value = "Asian - Black"
$ws.Cells.Item(26, 2).Value = -6.69
$ws.Cells.Item(26, 3).Value = 1.443
$ws.Cells.Item(26, 4).Value = 1669
$ws.Cells.Item(26, 5).Value = -4.636
$ws.Cells.Item(26, 6).Value = "'" + "3.7629e-05"
$ws.Cells.Item(26, 7).Value = "'" + "TRUE"
$ws.Cells.Item(26, 8).Value = "BiW_C"

$ws.Cells.Item(27, 1).Value = "Asian - LatinX"
$ws.Cells.Item(27, 2).Value = -5.635
$ws.Cells.Item(27, 3).Value = 1.853
$ws.Cells.Item(27, 4).Value = 1669
$ws.Cells.Item(27, 5).Value = -3.041
$ws.Cells.Item(27, 6).Value = "'" + "2.0224e-02"
$ws.Cells.Item(27, 7).Value = "'" + "TRUE"
$ws.Cells.Item(27, 8).Value = "BiW_C"

$ws.Cells.Item(28, 1).Value = "Asian - Other"
$ws.Cells.Item(28, 2).Value = -4.678
$ws.Cells.Item(28, 3).Value = 2.689
$ws.Cells.Item(28, 4).Value = 1669
$ws.Cells.Item(28, 5).Value = -1.74
$ws.Cells.Item(28, 6).Value = "'" + "4.0984e-01"
$ws.Cells.Item(28, 7).Value = "'" + "FALSE"
$ws.Cells.Item(28, 8).Value = "BiW_C"

$ws.Cells.Item(29, 1).Value = "Black - LatinX"
$ws.Cells.Item(29, 2).Value = 1.055
$ws.Cells.Item(29, 3).Value = 1.418
$ws.Cells.Item(29, 4).Value = 1669
$ws.Cells.Item(29, 5).Value = 0.744
$ws.Cells.Item(29, 6).Value = "'" + "9.4605e-01"
$ws.Cells.Item(29, 7).Value = "'" + "FALSE"
$ws.Cells.Item(29, 8).Value = "BiW_C"

$ws.Cells.Item(30, 1).Value = "Black - Other"
$ws.Cells.Item(30, 2).Value = 2.013
$ws.Cells.Item(30, 3).Value = 2.404
$ws.Cells.Item(30, 4).Value = 1669
$ws.Cells.Item(30, 5).Value = 0.837
$ws.Cells.Item(30, 6).Value = "'" + "9.1905e-01"
$ws.Cells.Item(30, 7).Value = "'" + "FALSE"
$ws.Cells.Item(30, 8).Value = "BiW_C"

$ws.Cells.Item(31, 1).Value = "LatinX - Other"
$ws.Cells.Item(31, 2).Value = 0.958
$ws.Cells.Item(31, 3).Value = 2.674
$ws.Cells.Item(31, 4).Value = 1669
$ws.Cells.Item(31, 5).Value = 0.358
$ws.Cells.Item(31, 6).Value = "'" + "9.9648e-01"
$ws.Cells.Item(31, 7).Value = "'" + "FALSE"
$ws.Cells.Item(31, 8).Value = "BiW_C"

$ws.Cells.Item(32, 1).Value = "White - Asian"
$ws.Cells.Item(32, 2).Value = 2.249
$ws.Cells.Item(32, 3).Value = 1.478
$ws.Cells.Item(32, 4).Value = 1669
$ws.Cells.Item(32, 5).Value = 1.522
$ws.Cells.Item(32, 6).Value = "'" + "5.4826e-01"
$ws.Cells.Item(32, 7).Value = "'" + "FALSE"
$ws.Cells.Item(32, 8).Value = "GoSub_C"

$ws.Cells.Item(33, 1).Value = "White - Black"
$ws.Cells.Item(33, 2).Value = -4.395
$ws.Cells.Item(33, 3).Value = 0.725
$ws.Cells.Item(33, 4).Value = 1669
$ws.Cells.Item(33, 5).Value = -6.066
$ws.Cells.Item(33, 6).Value = "'" + "1.6216e-08"
$ws.Cells.Item(33, 7).Value = "'" + "TRUE"
$ws.Cells.Item(33, 8).Value = "GoSub_C"

$ws.Cells.Item(34, 1).Value = "White - LatinX"
$ws.Cells.Item(34, 2).Value = 0.151
$ws.Cells.Item(34, 3).Value = 1.45
$ws.Cells.Item(34, 4).Value = 1669
$ws.Cells.Item(34, 5).Value = 0.104
$ws.Cells.Item(34, 6).Value = "'" + "9.9997e-01"
$ws.Cells.Item(34, 7).Value = "'" + "FALSE"
$ws.Cells.Item(34, 8).Value = "GoSub_C"

$ws.Cells.Item(35, 1).Value = "White - Other"
$ws.Cells.Item(35, 2).Value = 3.387
$ws.Cells.Item(35, 3).Value = 2.541
$ws.Cells.Item(35, 4).Value = 1669
$ws.Cells.Item(35, 5).Value = 1.333
$ws.Cells.Item(35, 6).Value = "'" + "6.7065e-01"
$ws.Cells.Item(35, 7).Value = "'" + "FALSE"
$ws.Cells.Item(35, 8).Value = "GoSub_C"

$ws.Cells.Item(36, 1).Value = "Asian - Black"
$ws.Cells.Item(36, 2).Value = -6.644
$ws.Cells.Item(36, 3).Value = 1.549
$ws.Cells.Item(36, 4).Value = 1669
$ws.Cells.Item(36, 5).Value = -4.289
$ws.Cells.Item(36, 6).Value = "'" + "1.8392e-04"
$ws.Cells.Item(36, 7).Value = "'" + "TRUE"
$ws.Cells.Item(36, 8).Value = "GoSub_C"

$ws.Cells.Item(37, 1).Value = "Asian - LatinX"
$ws.Cells.Item(37, 2).Value = -2.098
$ws.Cells.Item(37, 3).Value = 1.989
$ws.Cells.Item(37, 4).Value = 1669
$ws.Cells.Item(37, 5).Value = -1.055
$ws.Cells.Item(37, 6).Value = "'" + "8.2957e-01"
$ws.Cells.Item(37, 7).Value = "'" + "FALSE"
$ws.Cells.Item(37, 8).Value = "GoSub_C"

$ws.Cells.Item(38, 1).Value = "Asian - Other"
$ws.Cells.Item(38, 2).Value = 1.137
$ws.Cells.Item(38, 3).Value = 2.886
$ws.Cells.Item(38, 4).Value = 1669
$ws.Cells.Item(38, 5).Value = 0.394
$ws.Cells.Item(38, 6).Value = "'" + "9.9491e-01"
$ws.Cells.Item(38, 7).Value = "'" + "FALSE"
$ws.Cells.Item(38, 8).Value = "GoSub_C"

$ws.Cells.Item(39, 1).Value = "Black - LatinX"
$ws.Cells.Item(39, 2).Value = 4.546
$ws.Cells.Item(39, 3).Value = 1.522
$ws.Cells.Item(39, 4).Value = 1669
$ws.Cells.Item(39, 5).Value = 2.987
$ws.Cells.Item(39, 6).Value = "'" + "2.3847e-02"
$ws.Cells.Item(39, 7).Value = "'" + "TRUE"
$ws.Cells.Item(39, 8).Value = "GoSub_C"

$ws.Cells.Item(40, 1).Value = "Black - Other"
$ws.Cells.Item(40, 2).Value = 7.781
$ws.Cells.Item(40, 3).Value = 2.581
$ws.Cells.Item(40, 4).Value = 1669
$ws.Cells.Item(40, 5).Value = 3.015
$ws.Cells.Item(40, 6).Value = "'" + "2.1891e-02"
$ws.Cells.Item(40, 7).Value = "'" + "TRUE"
$ws.Cells.Item(40, 8).Value = "GoSub_C"

$ws.Cells.Item(41, 1).Value = "LatinX - Other"
$ws.Cells.Item(41, 2).Value = 3.235
$ws.Cells.Item(41, 3).Value = 2.87
$ws.Cells.Item(41, 4).Value = 1669
$ws.Cells.Item(41, 5).Value = 1.127
$ws.Cells.Item(41, 6).Value = "'" + "7.9217e-01"
$ws.Cells.Item(41, 7).Value = "'" + "FALSE"
$ws.Cells.Item(41, 8).Value = "GoSub_C"

$ws.Cells.Item(42, 1).Value = "White - Asian"
$ws.Cells.Item(42, 2).Value = 0.246
$ws.Cells.Item(42, 3).Value = 0.527
$ws.Cells.Item(42, 4).Value = 1669
$ws.Cells.Item(42, 5).Value = 0.467
$ws.Cells.Item(42, 6).Value = "'" + "9.9023e-01"
$ws.Cells.Item(42, 7).Value = "'" + "FALSE"
$ws.Cells.Item(42, 8).Value = "NRB_L"

$ws.Cells.Item(43, 1).Value = "White - Black"
$ws.Cells.Item(43, 2).Value = -2.43
$ws.Cells.Item(43, 3).Value = 0.258
$ws.Cells.Item(43, 4).Value = 1669
$ws.Cells.Item(43, 5).Value = -9.403
$ws.Cells.Item(43, 6).Value = "'" + "1.6522e-12"
$ws.Cells.Item(43, 7).Value = "'" + "TRUE"
$ws.Cells.Item(43, 8).Value = "NRB_L"

$ws.Cells.Item(44, 1).Value = "White - LatinX"
$ws.Cells.Item(44, 2).Value = -0.454
$ws.Cells.Item(44, 3).Value = 0.517
$ws.Cells.Item(44, 4).Value = 1669
$ws.Cells.Item(44, 5).Value = -0.878
$ws.Cells.Item(44, 6).Value = "'" + "9.0515e-01"
$ws.Cells.Item(44, 7).Value = "'" + "FALSE"
$ws.Cells.Item(44, 8).Value = "NRB_L"

$ws.Cells.Item(45, 1).Value = "White - Other"
$ws.Cells.Item(45, 2).Value = -1.895
$ws.Cells.Item(45, 3).Value = 0.906
$ws.Cells.Item(45, 4).Value = 1669
$ws.Cells.Item(45, 5).Value = -2.091
$ws.Cells.Item(45, 6).Value = "'" + "2.2445e-01"
$ws.Cells.Item(45, 7).Value = "'" + "FALSE"
$ws.Cells.Item(45, 8).Value = "NRB_L"

$ws.Cells.Item(46, 1).Value = "Asian - Black"
$ws.Cells.Item(46, 2).Value = -2.676
$ws.Cells.Item(46, 3).Value = 0.553
$ws.Cells.Item(46, 4).Value = 1669
$ws.Cells.Item(46, 5).Value = -4.844
$ws.Cells.Item(46, 6).Value = "'" + "1.3729e-05"
$ws.Cells.Item(46, 7).Value = "'" + "TRUE"
$ws.Cells.Item(46, 8).Value = "NRB_L"

$ws.Cells.Item(47, 1).Value = "Asian - LatinX"
$ws.Cells.Item(47, 2).Value = -0.7
$ws.Cells.Item(47, 3).Value = 0.709
$ws.Cells.Item(47, 4).Value = 1669
$ws.Cells.Item(47, 5).Value = -0.987
$ws.Cells.Item(47, 6).Value = "'" + "8.6129e-01"
$ws.Cells.Item(47, 7).Value = "'" + "FALSE"
$ws.Cells.Item(47, 8).Value = "NRB_L"

$ws.Cells.Item(48, 1).Value = "Asian - Other"
$ws.Cells.Item(48, 2).Value = -2.141
$ws.Cells.Item(48, 3).Value = 1.029
$ws.Cells.Item(48, 4).Value = 1669
$ws.Cells.Item(48, 5).Value = -2.08
$ws.Cells.Item(48, 6).Value = "'" + "2.2920e-01"
$ws.Cells.Item(48, 7).Value = "'" + "FALSE"
$ws.Cells.Item(48, 8).Value = "NRB_L"

$ws.Cells.Item(49, 1).Value = "Black - LatinX"
$ws.Cells.Item(49, 2).Value = 1.976
$ws.Cells.Item(49, 3).Value = 0.543
$ws.Cells.Item(49, 4).Value = 1669
$ws.Cells.Item(49, 5).Value = 3.641
$ws.Cells.Item(49, 6).Value = "'" + "2.5953e-03"
$ws.Cells.Item(49, 7).Value = "'" + "TRUE"
$ws.Cells.Item(49, 8).Value = "NRB_L"

$ws.Cells.Item(50, 1).Value = "Black - Other"
$ws.Cells.Item(50, 2).Value = 0.535
$ws.Cells.Item(50, 3).Value = 0.921
$ws.Cells.Item(50, 4).Value = 1669
$ws.Cells.Item(50, 5).Value = 0.581
$ws.Cells.Item(50, 6).Value = "'" + "9.7783e-01"
$ws.Cells.Item(50, 7).Value = "'" + "FALSE"
$ws.Cells.Item(50, 8).Value = "NRB_L"

$ws.Cells.Item(51, 1).Value = "LatinX - Other"
$ws.Cells.Item(51, 2).Value = -1.441
$ws.Cells.Item(51, 3).Value = 1.024
$ws.Cells.Item(51, 4).Value = 1669
$ws.Cells.Item(51, 5).Value = -1.408
$ws.Cells.Item(51, 6).Value = "'" + "6.2287e-01"
$ws.Cells.Item(51, 7).Value = "'" + "FALSE"
$ws.Cells.Item(51, 8).Value = "NRB_L"

$ws.Cells.Item(52, 1).Value = "White - Asian"
$ws.Cells.Item(52, 2).Value = 1.649
$ws.Cells.Item(52, 3).Value = 0.273
$ws.Cells.Item(52, 4).Value = 1669
$ws.Cells.Item(52, 5).Value = 6.049
$ws.Cells.Item(52, 6).Value = "'" + "1.7970e-08"
$ws.Cells.Item(52, 7).Value = "'" + "TRUE"
$ws.Cells.Item(52, 8).Value = "ProS_L"

$ws.Cells.Item(53, 1).Value = "White - Black"
$ws.Cells.Item(53, 2).Value = 1.672
$ws.Cells.Item(53, 3).Value = 0.134
$ws.Cells.Item(53, 4).Value = 1669
$ws.Cells.Item(53, 5).Value = 12.51
$ws.Cells.Item(53, 6).Value = "'" + "1.6044e-12"
$ws.Cells.Item(53, 7).Value = "'" + "TRUE"
$ws.Cells.Item(53, 8).Value = "ProS_L"

$ws.Cells.Item(54, 1).Value = "White - LatinX"
$ws.Cells.Item(54, 2).Value = 1.093
$ws.Cells.Item(54, 3).Value = 0.267
$ws.Cells.Item(54, 4).Value = 1669
$ws.Cells.Item(54, 5).Value = 4.088
$ws.Cells.Item(54, 6).Value = "'" + "4.3766e-04"
$ws.Cells.Item(54, 7).Value = "'" + "TRUE"
$ws.Cells.Item(54, 8).Value = "ProS_L"

$ws.Cells.Item(55, 1).Value = "White - Other"
$ws.Cells.Item(55, 2).Value = 2.113
$ws.Cells.Item(55, 3).Value = 0.469
$ws.Cells.Item(55, 4).Value = 1669
$ws.Cells.Item(55, 5).Value = 4.508
$ws.Cells.Item(55, 6).Value = "'" + "6.8428e-05"
$ws.Cells.Item(55, 7).Value = "'" + "TRUE"
$ws.Cells.Item(55, 8).Value = "ProS_L"

$ws.Cells.Item(56, 1).Value = "Asian - Black"
$ws.Cells.Item(56, 2).Value = 0.023
$ws.Cells.Item(56, 3).Value = 0.286
$ws.Cells.Item(56, 4).Value = 1669
$ws.Cells.Item(56, 5).Value = 0.081
$ws.Cells.Item(56, 6).Value = "'" + "9.9999e-01"
$ws.Cells.Item(56, 7).Value = "'" + "FALSE"
$ws.Cells.Item(56, 8).Value = "ProS_L"

$ws.Cells.Item(57, 1).Value = "Asian - LatinX"
$ws.Cells.Item(57, 2).Value = -0.556
$ws.Cells.Item(57, 3).Value = 0.367
$ws.Cells.Item(57, 4).Value = 1669
$ws.Cells.Item(57, 5).Value = -1.515
$ws.Cells.Item(57, 6).Value = "'" + "5.5305e-01"
$ws.Cells.Item(57, 7).Value = "'" + "FALSE"
$ws.Cells.Item(57, 8).Value = "ProS_L"

$ws.Cells.Item(58, 1).Value = "Asian - Other"
$ws.Cells.Item(58, 2).Value = 0.464
$ws.Cells.Item(58, 3).Value = 0.532
$ws.Cells.Item(58, 4).Value = 1669
$ws.Cells.Item(58, 5).Value = 0.872
$ws.Cells.Item(58, 6).Value = "'" + "9.0725e-01"
$ws.Cells.Item(58, 7).Value = "'" + "FALSE"
$ws.Cells.Item(58, 8).Value = "ProS_L"

$ws.Cells.Item(59, 1).Value = "Black - LatinX"
$ws.Cells.Item(59, 2).Value = -0.579
$ws.Cells.Item(59, 3).Value = 0.281
$ws.Cells.Item(59, 4).Value = 1669
$ws.Cells.Item(59, 5).Value = -2.062
$ws.Cells.Item(59, 6).Value = "'" + "2.3735e-01"
$ws.Cells.Item(59, 7).Value = "'" + "FALSE"
$ws.Cells.Item(59, 8).Value = "ProS_L"

$ws.Cells.Item(60, 1).Value = "Black - Other"
$ws.Cells.Item(60, 2).Value = 0.441
$ws.Cells.Item(60, 3).Value = 0.476
$ws.Cells.Item(60, 4).Value = 1669
$ws.Cells.Item(60, 5).Value = 0.926
$ws.Cells.Item(60, 6).Value = "'" + "8.8671e-01"
$ws.Cells.Item(60, 7).Value = "'" + "FALSE"
$ws.Cells.Item(60, 8).Value = "ProS_L"

$ws.Cells.Item(61, 1).Value = "LatinX - Other"
$ws.Cells.Item(61, 2).Value = 1.02
$ws.Cells.Item(61, 3).Value = 0.529
$ws.Cells.Item(61, 4).Value = 1669
$ws.Cells.Item(61, 5).Value = 1.926
$ws.Cells.Item(61, 6).Value = "'" + "3.0373e-01"
$ws.Cells.Item(61, 7).Value = "'" + "FALSE"
$ws.Cells.Item(61, 8).Value = "ProS_L"

$ws.Cells.Item(62, 1).Value = "White - Asian"
$ws.Cells.Item(62, 2).Value = 1.768
$ws.Cells.Item(62, 3).Value = 0.441
$ws.Cells.Item(62, 4).Value = 1669
$ws.Cells.Item(62, 5).Value = 4.013
$ws.Cells.Item(62, 6).Value = "'" + "5.9902e-04"
$ws.Cells.Item(62, 7).Value = "'" + "TRUE"
$ws.Cells.Item(62, 8).Value = "SelP_L"

$ws.Cells.Item(63, 1).Value = "White - Black"
$ws.Cells.Item(63, 2).Value = 1.944
$ws.Cells.Item(63, 3).Value = 0.216
$ws.Cells.Item(63, 4).Value = 1669
$ws.Cells.Item(63, 5).Value = 8.998
$ws.Cells.Item(63, 6).Value = "'" + "1.6489e-12"
$ws.Cells.Item(63, 7).Value = "'" + "TRUE"
$ws.Cells.Item(63, 8).Value = "SelP_L"

$ws.Cells.Item(64, 1).Value = "White - LatinX"
$ws.Cells.Item(64, 2).Value = 0.184
$ws.Cells.Item(64, 3).Value = 0.432
$ws.Cells.Item(64, 4).Value = 1669
$ws.Cells.Item(64, 5).Value = 0.425
$ws.Cells.Item(64, 6).Value = "'" + "9.9320e-01"
$ws.Cells.Item(64, 7).Value = "'" + "FALSE"
$ws.Cells.Item(64, 8).Value = "SelP_L"

$ws.Cells.Item(65, 1).Value = "White - Other"
$ws.Cells.Item(65, 2).Value = 0.045
$ws.Cells.Item(65, 3).Value = 0.758
$ws.Cells.Item(65, 4).Value = 1669
$ws.Cells.Item(65, 5).Value = 0.06
$ws.Cells.Item(65, 6).Value = "'" + "1.0000e+00"
$ws.Cells.Item(65, 7).Value = "'" + "FALSE"
$ws.Cells.Item(65, 8).Value = "SelP_L"

$ws.Cells.Item(66, 1).Value = "Asian - Black"
$ws.Cells.Item(66, 2).Value = 0.176
$ws.Cells.Item(66, 3).Value = 0.462
$ws.Cells.Item(66, 4).Value = 1669
$ws.Cells.Item(66, 5).Value = 0.381
$ws.Cells.Item(66, 6).Value = "'" + "9.9555e-01"
$ws.Cells.Item(66, 7).Value = "'" + "FALSE"
$ws.Cells.Item(66, 8).Value = "SelP_L"

$ws.Cells.Item(67, 1).Value = "Asian - LatinX"
$ws.Cells.Item(67, 2).Value = -1.584
$ws.Cells.Item(67, 3).Value = 0.593
$ws.Cells.Item(67, 4).Value = 1669
$ws.Cells.Item(67, 5).Value = -2.672
$ws.Cells.Item(67, 6).Value = "'" + "5.8637e-02"
$ws.Cells.Item(67, 7).Value = "'" + "FALSE"
$ws.Cells.Item(67, 8).Value = "SelP_L"

$ws.Cells.Item(68, 1).Value = "Asian - Other"
$ws.Cells.Item(68, 2).Value = -1.723
$ws.Cells.Item(68, 3).Value = 0.86
$ws.Cells.Item(68, 4).Value = 1669
$ws.Cells.Item(68, 5).Value = -2.002
$ws.Cells.Item(68, 6).Value = "'" + "2.6528e-01"
$ws.Cells.Item(68, 7).Value = "'" + "FALSE"
$ws.Cells.Item(68, 8).Value = "SelP_L"

$ws.Cells.Item(69, 1).Value = "Black - LatinX"
$ws.Cells.Item(69, 2).Value = -1.76
$ws.Cells.Item(69, 3).Value = 0.454
$ws.Cells.Item(69, 4).Value = 1669
$ws.Cells.Item(69, 5).Value = -3.879
$ws.Cells.Item(69, 6).Value = "'" + "1.0301e-03"
$ws.Cells.Item(69, 7).Value = "'" + "TRUE"
$ws.Cells.Item(69, 8).Value = "SelP_L"

$ws.Cells.Item(70, 1).Value = "Black - Other"
$ws.Cells.Item(70, 2).Value = -1.899
$ws.Cells.Item(70, 3).Value = 0.769
$ws.Cells.Item(70, 4).Value = 1669
$ws.Cells.Item(70, 5).Value = -2.468
$ws.Cells.Item(70, 6).Value = "'" + "9.8552e-02"
$ws.Cells.Item(70, 7).Value = "'" + "FALSE"
$ws.Cells.Item(70, 8).Value = "SelP_L"

$ws.Cells.Item(71, 1).Value = "LatinX - Other"
$ws.Cells.Item(71, 2).Value = -0.138
$ws.Cells.Item(71, 3).Value = 0.856
$ws.Cells.Item(71, 4).Value = 1669
$ws.Cells.Item(71, 5).Value = -0.162
$ws.Cells.Item(71, 6).Value = "'" + "9.9985e-01"
$ws.Cells.Item(71, 7).Value = "'" + "FALSE"
$ws.Cells.Item(71, 8).Value = "SelP_L"

$ws.Cells.Item(72, 1).Value = "White - Asian"
$ws.Cells.Item(72, 2).Value = 4.764
$ws.Cells.Item(72, 3).Value = 0.96
$ws.Cells.Item(72, 4).Value = 1669
$ws.Cells.Item(72, 5).Value = 4.965
$ws.Cells.Item(72, 6).Value = "'" + "7.4964e-06"
$ws.Cells.Item(72, 7).Value = "'" + "TRUE"
$ws.Cells.Item(72, 8).Value = "SelM_L"

$ws.Cells.Item(73, 1).Value = "White - Black"
$ws.Cells.Item(73, 2).Value = -2.678
$ws.Cells.Item(73, 3).Value = 0.47
$ws.Cells.Item(73, 4).Value = 1669
$ws.Cells.Item(73, 5).Value = -5.692
$ws.Cells.Item(73, 6).Value = "'" + "1.4767e-07"
$ws.Cells.Item(73, 7).Value = "'" + "TRUE"
$ws.Cells.Item(73, 8).Value = "SelM_L"

$ws.Cells.Item(74, 1).Value = "White - LatinX"
$ws.Cells.Item(74, 2).Value = 1.059
$ws.Cells.Item(74, 3).Value = 0.941
$ws.Cells.Item(74, 4).Value = 1669
$ws.Cells.Item(74, 5).Value = 1.126
$ws.Cells.Item(74, 6).Value = "'" + "7.9306e-01"
$ws.Cells.Item(74, 7).Value = "'" + "FALSE"
$ws.Cells.Item(74, 8).Value = "SelM_L"

$ws.Cells.Item(75, 1).Value = "White - Other"
$ws.Cells.Item(75, 2).Value = 0.607
$ws.Cells.Item(75, 3).Value = 1.65
$ws.Cells.Item(75, 4).Value = 1669
$ws.Cells.Item(75, 5).Value = 0.368
$ws.Cells.Item(75, 6).Value = "'" + "9.9609e-01"
$ws.Cells.Item(75, 7).Value = "'" + "FALSE"
$ws.Cells.Item(75, 8).Value = "SelM_L"

$ws.Cells.Item(76, 1).Value = "Asian - Black"
$ws.Cells.Item(76, 2).Value = -7.442
$ws.Cells.Item(76, 3).Value = 1.006
$ws.Cells.Item(76, 4).Value = 1669
$ws.Cells.Item(76, 5).Value = -7.4
$ws.Cells.Item(76, 6).Value = "'" + "3.7907e-12"
$ws.Cells.Item(76, 7).Value = "'" + "TRUE"
$ws.Cells.Item(76, 8).Value = "SelM_L"

$ws.Cells.Item(77, 1).Value = "Asian - LatinX"
$ws.Cells.Item(77, 2).Value = -3.705
$ws.Cells.Item(77, 3).Value = 1.291
$ws.Cells.Item(77, 4).Value = 1669
$ws.Cells.Item(77, 5).Value = -2.869
$ws.Cells.Item(77, 6).Value = "'" + "3.3894e-02"
$ws.Cells.Item(77, 7).Value = "'" + "TRUE"
$ws.Cells.Item(77, 8).Value = "SelM_L"

$ws.Cells.Item(78, 1).Value = "Asian - Other"
$ws.Cells.Item(78, 2).Value = -4.157
$ws.Cells.Item(78, 3).Value = 1.874
$ws.Cells.Item(78, 4).Value = 1669
$ws.Cells.Item(78, 5).Value = -2.218
$ws.Cells.Item(78, 6).Value = "'" + "1.7330e-01"
$ws.Cells.Item(78, 7).Value = "'" + "FALSE"
$ws.Cells.Item(78, 8).Value = "SelM_L"

$ws.Cells.Item(79, 1).Value = "Black - LatinX"
$ws.Cells.Item(79, 2).Value = 3.737
$ws.Cells.Item(79, 3).Value = 0.988
$ws.Cells.Item(79, 4).Value = 1669
$ws.Cells.Item(79, 5).Value = 3.782
$ws.Cells.Item(79, 6).Value = "'" + "1.5110e-03"
$ws.Cells.Item(79, 7).Value = "'" + "TRUE"
$ws.Cells.Item(79, 8).Value = "SelM_L"

$ws.Cells.Item(80, 1).Value = "Black - Other"
$ws.Cells.Item(80, 2).Value = 3.285
$ws.Cells.Item(80, 3).Value = 1.676
$ws.Cells.Item(80, 4).Value = 1669
$ws.Cells.Item(80, 5).Value = 1.96
$ws.Cells.Item(80, 6).Value = "'" + "2.8601e-01"
$ws.Cells.Item(80, 7).Value = "'" + "FALSE"
$ws.Cells.Item(80, 8).Value = "SelM_L"

$ws.Cells.Item(81, 1).Value = "LatinX - Other"
$ws.Cells.Item(81, 2).Value = -0.452
$ws.Cells.Item(81, 3).Value = 1.864
$ws.Cells.Item(81, 4).Value = 1669
$ws.Cells.Item(81, 5).Value = -0.243
$ws.Cells.Item(81, 6).Value = "'" + "9.9923e-01"
$ws.Cells.Item(81, 7).Value = "'" + "FALSE"
$ws.Cells.Item(81, 8).Value = "SelM_L"

$ws.Cells.Item(82, 1).Value = "White - Asian"
$ws.Cells.Item(82, 2).Value = 1.988
$ws.Cells.Item(82, 3).Value = 1.095
$ws.Cells.Item(82, 4).Value = 1669
$ws.Cells.Item(82, 5).Value = 1.816
$ws.Cells.Item(82, 6).Value = "'" + "3.6456e-01"
$ws.Cells.Item(82, 7).Value = "'" + "FALSE"
$ws.Cells.Item(82, 8).Value = "SnasM_C"

$ws.Cells.Item(83, 1).Value = "White - Black"
$ws.Cells.Item(83, 2).Value = -8.046
$ws.Cells.Item(83, 3).Value = 0.537
$ws.Cells.Item(83, 4).Value = 1669
$ws.Cells.Item(83, 5).Value = -14.991
$ws.Cells.Item(83, 6).Value = "'" + "1.6044e-12"
$ws.Cells.Item(83, 7).Value = "'" + "TRUE"
$ws.Cells.Item(83, 8).Value = "SnasM_C"

$ws.Cells.Item(84, 1).Value = "White - LatinX"
$ws.Cells.Item(84, 2).Value = -0.051
$ws.Cells.Item(84, 3).Value = 1.074
$ws.Cells.Item(84, 4).Value = 1669
$ws.Cells.Item(84, 5).Value = -0.048
$ws.Cells.Item(84, 6).Value = "'" + "1.0000e+00"
$ws.Cells.Item(84, 7).Value = "'" + "FALSE"
$ws.Cells.Item(84, 8).Value = "SnasM_C"

$ws.Cells.Item(85, 1).Value = "White - Other"
$ws.Cells.Item(85, 2).Value = -0.873
$ws.Cells.Item(85, 3).Value = 1.882
$ws.Cells.Item(85, 4).Value = 1669
$ws.Cells.Item(85, 5).Value = -0.464
$ws.Cells.Item(85, 6).Value = "'" + "9.9048e-01"
$ws.Cells.Item(85, 7).Value = "'" + "FALSE"
$ws.Cells.Item(85, 8).Value = "SnasM_C"

$ws.Cells.Item(86, 1).Value = "Asian - Black"
$ws.Cells.Item(86, 2).Value = -10.034
$ws.Cells.Item(86, 3).Value = 1.147
$ws.Cells.Item(86, 4).Value = 1669
$ws.Cells.Item(86, 5).Value = -8.745
$ws.Cells.Item(86, 6).Value = "'" + "1.6447e-12"
$ws.Cells.Item(86, 7).Value = "'" + "TRUE"
$ws.Cells.Item(86, 8).Value = "SnasM_C"

$ws.Cells.Item(87, 1).Value = "Asian - LatinX"
$ws.Cells.Item(87, 2).Value = -2.039
$ws.Cells.Item(87, 3).Value = 1.473
$ws.Cells.Item(87, 4).Value = 1669
$ws.Cells.Item(87, 5).Value = -1.384
$ws.Cells.Item(87, 6).Value = "'" + "6.3797e-01"
$ws.Cells.Item(87, 7).Value = "'" + "FALSE"
$ws.Cells.Item(87, 8).Value = "SnasM_C"

$ws.Cells.Item(88, 1).Value = "Asian - Other"
$ws.Cells.Item(88, 2).Value = -2.862
$ws.Cells.Item(88, 3).Value = 2.138
$ws.Cells.Item(88, 4).Value = 1669
$ws.Cells.Item(88, 5).Value = -1.338
$ws.Cells.Item(88, 6).Value = "'" + "6.6716e-01"
$ws.Cells.Item(88, 7).Value = "'" + "FALSE"
$ws.Cells.Item(88, 8).Value = "SnasM_C"

$ws.Cells.Item(89, 1).Value = "Black - LatinX"
$ws.Cells.Item(89, 2).Value = 7.995
$ws.Cells.Item(89, 3).Value = 1.127
$ws.Cells.Item(89, 4).Value = 1669
$ws.Cells.Item(89, 5).Value = 7.092
$ws.Cells.Item(89, 6).Value = "'" + "2.1127e-11"
$ws.Cells.Item(89, 7).Value = "'" + "TRUE"
$ws.Cells.Item(89, 8).Value = "SnasM_C"

$ws.Cells.Item(90, 1).Value = "Black - Other"
$ws.Cells.Item(90, 2).Value = 7.173
$ws.Cells.Item(90, 3).Value = 1.912
$ws.Cells.Item(90, 4).Value = 1669
$ws.Cells.Item(90, 5).Value = 3.752
$ws.Cells.Item(90, 6).Value = "'" + "1.6979e-03"
$ws.Cells.Item(90, 7).Value = "'" + "TRUE"
$ws.Cells.Item(90, 8).Value = "SnasM_C"

$ws.Cells.Item(91, 1).Value = "LatinX - Other"
$ws.Cells.Item(91, 2).Value = -0.822
$ws.Cells.Item(91, 3).Value = 2.126
$ws.Cells.Item(91, 4).Value = 1669
$ws.Cells.Item(91, 5).Value = -0.387
$ws.Cells.Item(91, 6).Value = "'" + "9.9527e-01"
$ws.Cells.Item(91, 7).Value = "'" + "FALSE"
$ws.Cells.Item(91, 8).Value = "SnasM_C"

$ws.Cells.Item(92, 1).Value = "White - Asian"
$ws.Cells.Item(92, 2).Value = 2.284
$ws.Cells.Item(92, 3).Value = 1.277
$ws.Cells.Item(92, 4).Value = 1669
$ws.Cells.Item(92, 5).Value = 1.789
$ws.Cells.Item(92, 6).Value = "'" + "3.8016e-01"
$ws.Cells.Item(92, 7).Value = "'" + "FALSE"
$ws.Cells.Item(92, 8).Value = "TrSman_C"

$ws.Cells.Item(93, 1).Value = "White - Black"
$ws.Cells.Item(93, 2).Value = -4.607
$ws.Cells.Item(93, 3).Value = 0.626
$ws.Cells.Item(93, 4).Value = 1669
$ws.Cells.Item(93, 5).Value = -7.36
$ws.Cells.Item(93, 6).Value = "'" + "4.5137e-12"
$ws.Cells.Item(93, 7).Value = "'" + "TRUE"
$ws.Cells.Item(93, 8).Value = "TrSman_C"

$ws.Cells.Item(94, 1).Value = "White - LatinX"
$ws.Cells.Item(94, 2).Value = 0.608
$ws.Cells.Item(94, 3).Value = 1.252
$ws.Cells.Item(94, 4).Value = 1669
$ws.Cells.Item(94, 5).Value = 0.485
$ws.Cells.Item(94, 6).Value = "'" + "9.8872e-01"
$ws.Cells.Item(94, 7).Value = "'" + "FALSE"
$ws.Cells.Item(94, 8).Value = "TrSman_C"

$ws.Cells.Item(95, 1).Value = "White - Other"
$ws.Cells.Item(95, 2).Value = 2.398
$ws.Cells.Item(95, 3).Value = 2.195
$ws.Cells.Item(95, 4).Value = 1669
$ws.Cells.Item(95, 5).Value = 1.092
$ws.Cells.Item(95, 6).Value = "'" + "8.1055e-01"
$ws.Cells.Item(95, 7).Value = "'" + "FALSE"
$ws.Cells.Item(95, 8).Value = "TrSman_C"

$ws.Cells.Item(96, 1).Value = "Asian - Black"
$ws.Cells.Item(96, 2).Value = -6.891
$ws.Cells.Item(96, 3).Value = 1.338
$ws.Cells.Item(96, 4).Value = 1669
$ws.Cells.Item(96, 5).Value = -5.15
$ws.Cells.Item(96, 6).Value = "'" + "2.8966e-06"
$ws.Cells.Item(96, 7).Value = "'" + "TRUE"
$ws.Cells.Item(96, 8).Value = "TrSman_C"

$ws.Cells.Item(97, 1).Value = "Asian - LatinX"
$ws.Cells.Item(97, 2).Value = -1.677
$ws.Cells.Item(97, 3).Value = 1.718
$ws.Cells.Item(97, 4).Value = 1669
$ws.Cells.Item(97, 5).Value = -0.976
$ws.Cells.Item(97, 6).Value = "'" + "8.6621e-01"
$ws.Cells.Item(97, 7).Value = "'" + "FALSE"
$ws.Cells.Item(97, 8).Value = "TrSman_C"

$ws.Cells.Item(98, 1).Value = "Asian - Other"
$ws.Cells.Item(98, 2).Value = 0.113
$ws.Cells.Item(98, 3).Value = 2.493
$ws.Cells.Item(98, 4).Value = 1669
$ws.Cells.Item(98, 5).Value = 0.046
$ws.Cells.Item(98, 6).Value = "'" + "1.0000e+00"
$ws.Cells.Item(98, 7).Value = "'" + "FALSE"
$ws.Cells.Item(98, 8).Value = "TrSman_C"

$ws.Cells.Item(99, 1).Value = "Black - LatinX"
$ws.Cells.Item(99, 2).Value = 5.214
$ws.Cells.Item(99, 3).Value = 1.315
$ws.Cells.Item(99, 4).Value = 1669
$ws.Cells.Item(99, 5).Value = 3.966
$ws.Cells.Item(99, 6).Value = "'" + "7.2486e-04"
$ws.Cells.Item(99, 7).Value = "'" + "TRUE"
$ws.Cells.Item(99, 8).Value = "TrSman_C"

$ws.Cells.Item(100, 1).Value = "Black - Other"
$ws.Cells.Item(100, 2).Value = 7.004
$ws.Cells.Item(100, 3).Value = 2.229
$ws.Cells.Item(100, 4).Value = 1669
$ws.Cells.Item(100, 5).Value = 3.142
$ws.Cells.Item(100, 6).Value = "'" + "1.4717e-02"
$ws.Cells.Item(100, 7).Value = "'" + "TRUE"
$ws.Cells.Item(100, 8).Value = "TrSman_C"

$ws.Cells.Item(101, 1).Value = "LatinX - Other"
$ws.Cells.Item(101, 2).Value = 1.79
$ws.Cells.Item(101, 3).Value = 2.48
$ws.Cells.Item(101, 4).Value = 1669
$ws.Cells.Item(101, 5).Value = 0.722
$ws.Cells.Item(101, 6).Value = "'" + "9.5149e-01"
$ws.Cells.Item(101, 7).Value = "'" + "FALSE"
$ws.Cells.Item(101, 8).Value = "TrSman_C"

$ws.Cells.Item(102, 1).Value = "White - Asian"
$ws.Cells.Item(102, 2).Value = 0.882
$ws.Cells.Item(102, 3).Value = 1.354
$ws.Cells.Item(102, 4).Value = 1669
$ws.Cells.Item(102, 5).Value = 0.651
$ws.Cells.Item(102, 6).Value = "'" + "9.6644e-01"
$ws.Cells.Item(102, 7).Value = "'" + "FALSE"
$ws.Cells.Item(102, 8).Value = "TrTr_C"

$ws.Cells.Item(103, 1).Value = "White - Black"
$ws.Cells.Item(103, 2).Value = -7.139
$ws.Cells.Item(103, 3).Value = 0.664
$ws.Cells.Item(103, 4).Value = 1669
$ws.Cells.Item(103, 5).Value = -10.754
$ws.Cells.Item(103, 6).Value = "'" + "1.6478e-12"
$ws.Cells.Item(103, 7).Value = "'" + "TRUE"
$ws.Cells.Item(103, 8).Value = "TrTr_C"

$ws.Cells.Item(104, 1).Value = "White - LatinX"
$ws.Cells.Item(104, 2).Value = 0.455
$ws.Cells.Item(104, 3).Value = 1.328
$ws.Cells.Item(104, 4).Value = 1669
$ws.Cells.Item(104, 5).Value = 0.342
$ws.Cells.Item(104, 6).Value = "'" + "9.9706e-01"
$ws.Cells.Item(104, 7).Value = "'" + "FALSE"
$ws.Cells.Item(104, 8).Value = "TrTr_C"

$ws.Cells.Item(105, 1).Value = "White - Other"
$ws.Cells.Item(105, 2).Value = 1.592
$ws.Cells.Item(105, 3).Value = 2.328
$ws.Cells.Item(105, 4).Value = 1669
$ws.Cells.Item(105, 5).Value = 0.684
$ws.Cells.Item(105, 6).Value = "'" + "9.5999e-01"
$ws.Cells.Item(105, 7).Value = "'" + "FALSE"
$ws.Cells.Item(105, 8).Value = "TrTr_C"

$ws.Cells.Item(106, 1).Value = "Asian - Black"
$ws.Cells.Item(106, 2).Value = -8.021
$ws.Cells.Item(106, 3).Value = 1.419
$ws.Cells.Item(106, 4).Value = 1669
$ws.Cells.Item(106, 5).Value = -5.652
$ws.Cells.Item(106, 6).Value = "'" + "1.8603e-07"
$ws.Cells.Item(106, 7).Value = "'" + "TRUE"
$ws.Cells.Item(106, 8).Value = "TrTr_C"

$ws.Cells.Item(107, 1).Value = "Asian - LatinX"
$ws.Cells.Item(107, 2).Value = -0.427
$ws.Cells.Item(107, 3).Value = 1.822
$ws.Cells.Item(107, 4).Value = 1669
$ws.Cells.Item(107, 5).Value = -0.234
$ws.Cells.Item(107, 6).Value = "'" + "9.9933e-01"
$ws.Cells.Item(107, 7).Value = "'" + "FALSE"
$ws.Cells.Item(107, 8).Value = "TrTr_C"

$ws.Cells.Item(108, 1).Value = "Asian - Other"
$ws.Cells.Item(108, 2).Value = 0.71
$ws.Cells.Item(108, 3).Value = 2.644
$ws.Cells.Item(108, 4).Value = 1669
$ws.Cells.Item(108, 5).Value = 0.269
$ws.Cells.Item(108, 6).Value = "'" + "9.9886e-01"
$ws.Cells.Item(108, 7).Value = "'" + "FALSE"
$ws.Cells.Item(108, 8).Value = "TrTr_C"

$ws.Cells.Item(109, 1).Value = "Black - LatinX"
$ws.Cells.Item(109, 2).Value = 7.594
$ws.Cells.Item(109, 3).Value = 1.394
$ws.Cells.Item(109, 4).Value = 1669
$ws.Cells.Item(109, 5).Value = 5.446
$ws.Cells.Item(109, 6).Value = "'" + "5.8989e-07"
$ws.Cells.Item(109, 7).Value = "'" + "TRUE"
$ws.Cells.Item(109, 8).Value = "TrTr_C"

$ws.Cells.Item(110, 1).Value = "Black - Other"
$ws.Cells.Item(110, 2).Value = 8.731
$ws.Cells.Item(110, 3).Value = 2.365
$ws.Cells.Item(110, 4).Value = 1669
$ws.Cells.Item(110, 5).Value = 3.693
$ws.Cells.Item(110, 6).Value = "'" + "2.1323e-03"
$ws.Cells.Item(110, 7).Value = "'" + "TRUE"
$ws.Cells.Item(110, 8).Value = "TrTr_C"

$ws.Cells.Item(111, 1).Value = "LatinX - Other"
$ws.Cells.Item(111, 2).Value = 1.138
$ws.Cells.Item(111, 3).Value = 2.63
$ws.Cells.Item(111, 4).Value = 1669
$ws.Cells.Item(111, 5).Value = 0.433
$ws.Cells.Item(111, 6).Value = "'" + "9.9272e-01"
$ws.Cells.Item(111, 7).Value = "'" + "FALSE"
$ws.Cells.Item(111, 8).Value = "TrTr_C"

$ws.Cells.Item(112, 1).Value = "White - Asian"
$ws.Cells.Item(112, 2).Value = -1.146
$ws.Cells.Item(112, 3).Value = 0.671
$ws.Cells.Item(112, 4).Value = 1669
$ws.Cells.Item(112, 5).Value = -1.708
$ws.Cells.Item(112, 6).Value = "'" + "4.2933e-01"
$ws.Cells.Item(112, 7).Value = "'" + "FALSE"
$ws.Cells.Item(112, 8).Value = "TrTr_L"

$ws.Cells.Item(113, 1).Value = "White - Black"
$ws.Cells.Item(113, 2).Value = -2.658
$ws.Cells.Item(113, 3).Value = 0.329
$ws.Cells.Item(113, 4).Value = 1669
$ws.Cells.Item(113, 5).Value = -8.082
$ws.Cells.Item(113, 6).Value = "'" + "1.6780e-12"
$ws.Cells.Item(113, 7).Value = "'" + "TRUE"
$ws.Cells.Item(113, 8).Value = "TrTr_L"

$ws.Cells.Item(114, 1).Value = "White - LatinX"
$ws.Cells.Item(114, 2).Value = -0.582
$ws.Cells.Item(114, 3).Value = 0.658
$ws.Cells.Item(114, 4).Value = 1669
$ws.Cells.Item(114, 5).Value = -0.884
$ws.Cells.Item(114, 6).Value = "'" + "9.0285e-01"
$ws.Cells.Item(114, 7).Value = "'" + "FALSE"
$ws.Cells.Item(114, 8).Value = "TrTr_L"

$ws.Cells.Item(115, 1).Value = "White - Other"
$ws.Cells.Item(115, 2).Value = -0.551
$ws.Cells.Item(115, 3).Value = 1.153
$ws.Cells.Item(115, 4).Value = 1669
$ws.Cells.Item(115, 5).Value = -0.478
$ws.Cells.Item(115, 6).Value = "'" + "9.8938e-01"
$ws.Cells.Item(115, 7).Value = "'" + "FALSE"
$ws.Cells.Item(115, 8).Value = "TrTr_L"

$ws.Cells.Item(116, 1).Value = "Asian - Black"
$ws.Cells.Item(116, 2).Value = -1.513
$ws.Cells.Item(116, 3).Value = 0.703
$ws.Cells.Item(116, 4).Value = 1669
$ws.Cells.Item(116, 5).Value = -2.151
$ws.Cells.Item(116, 6).Value = "'" + "1.9907e-01"
$ws.Cells.Item(116, 7).Value = "'" + "FALSE"
$ws.Cells.Item(116, 8).Value = "TrTr_L"

$ws.Cells.Item(117, 1).Value = "Asian - LatinX"
$ws.Cells.Item(117, 2).Value = 0.564
$ws.Cells.Item(117, 3).Value = 0.903
$ws.Cells.Item(117, 4).Value = 1669
$ws.Cells.Item(117, 5).Value = 0.625
$ws.Cells.Item(117, 6).Value = "'" + "9.7116e-01"
$ws.Cells.Item(117, 7).Value = "'" + "FALSE"
$ws.Cells.Item(117, 8).Value = "TrTr_L"

$ws.Cells.Item(118, 1).Value = "Asian - Other"
$ws.Cells.Item(118, 2).Value = 0.595
$ws.Cells.Item(118, 3).Value = 1.31
$ws.Cells.Item(118, 4).Value = 1669
$ws.Cells.Item(118, 5).Value = 0.454
$ws.Cells.Item(118, 6).Value = "'" + "9.9124e-01"
$ws.Cells.Item(118, 7).Value = "'" + "FALSE"
$ws.Cells.Item(118, 8).Value = "TrTr_L"

$ws.Cells.Item(119, 1).Value = "Black - LatinX"
$ws.Cells.Item(119, 2).Value = 2.077
$ws.Cells.Item(119, 3).Value = 0.691
$ws.Cells.Item(119, 4).Value = 1669
$ws.Cells.Item(119, 5).Value = 3.006
$ws.Cells.Item(119, 6).Value = "'" + "2.2532e-02"
$ws.Cells.Item(119, 7).Value = "'" + "TRUE"
$ws.Cells.Item(119, 8).Value = "TrTr_L"

$ws.Cells.Item(120, 1).Value = "Black - Other"
$ws.Cells.Item(120, 2).Value = 2.107
$ws.Cells.Item(120, 3).Value = 1.171
$ws.Cells.Item(120, 4).Value = 1669
$ws.Cells.Item(120, 5).Value = 1.799
$ws.Cells.Item(120, 6).Value = "'" + "3.7442e-01"
$ws.Cells.Item(120, 7).Value = "'" + "FALSE"
$ws.Cells.Item(120, 8).Value = "TrTr_L"

$ws.Cells.Item(121, 1).Value = "LatinX - Other"
$ws.Cells.Item(121, 2).Value = 0.031
$ws.Cells.Item(121, 3).Value = 1.303
$ws.Cells.Item(121, 4).Value = 1669
$ws.Cells.Item(121, 5).Value = 0.024
$ws.Cells.Item(121, 6).Value = "'" + "1.0000e+00"
$ws.Cells.Item(121, 7).Value = "'" + "FALSE"
$ws.Cells.Item(121, 8).Value = "TrTr_L"
